$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (ownTeam, oppTeam) before the existing "batsman" column (D),
# shifting batsman..sr from D:I to F:K.
$ws.Columns("D:E").Insert()

# Insert a new row 2 (new most-recent match) above the existing data row,
# shifting the existing data row from row 2 to row 3.
$ws.Rows("2:2").Insert()

# The numeric-looking stat columns (totalRuns, totalBalls, total4s, total6s, sr)
# must stay stored as text, like the rest of the sheet - force Text format so
# Excel doesn't silently coerce these into real numbers.
$ws.Range("G2:K4").NumberFormat = "@"

# New header cells for the two inserted columns.
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New row 2: October 09 2020 match vs Delhi Capitals.
$ws.Range("A2").Value = " Sharjah"
$ws.Range("B2").Value = " October 09 2020"
$ws.Range("C2").Value = "Capitals won by 46 runs"
$ws.Range("D2").Value = "Rajasthan Royals"
$ws.Range("E2").Value = "Delhi Capitals"
$ws.Range("F2").Value = "Yashasvi Jaiswal "
$ws.Range("G2").Value = "34"
$ws.Range("H2").Value = "36"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "2"
$ws.Range("K2").Value = "94.44"

# Existing row (now row 3) gains ownTeam / oppTeam values.
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "Chennai Super Kings"

# New row 4: October 06 2020 match vs Mumbai Indians.
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 06 2020"
$ws.Range("C4").Value = "Mumbai won by 57 runs"
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "Yashasvi Jaiswal "
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "2"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "0.00"
